# Actualización automática 2025-05-30 16:20:08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Column I width: 9 -> 11 (raw OOXML width units) ---
# Excel's ColumnWidth property (characters) differs from the stored <col width>
# by the standard ~0.8333 padding offset at the default font, so asking for
# 10.1667 yields a stored width of exactly 11.
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666

# --- Row 5 ---
$ws.Range("E5").Value = 217.58
$ws.Range("L5").Value = 1149.35
$ws.Range("N5").Value = 547.43

# --- Row 6 ---
$ws.Range("E6").Value = 69.45
$ws.Range("F6").Value = 52.25
$ws.Range("I6").Value = 15.68
$ws.Range("K6").Value = 2919.57
$ws.Range("L6").Value = 3782.45

# --- Row 7 ---
$ws.Range("L7").Value = 48.37
$ws.Range("N7").Value = 547.43

# --- Row 8 ---
$ws.Range("C8").Value = 1166.4
$ws.Range("D8").Value = 6044.16
$ws.Range("J8").Value = -199.99

# --- Row 10 ---
$ws.Range("D10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("L10").Value = 9002.94
$ws.Range("N10").Value = 2189.72

# --- Row 12 ---
$ws.Range("L12").Value = 1565.15

# --- Row 13 ---
$ws.Range("I13").Value = 57.46
$ws.Range("K13").Value = 3025.16
$ws.Range("L13").Value = 5769.32
$ws.Range("M13").Value = 182.61

# --- Row 15 ---
$ws.Range("L15").Value = 288.17

# --- Row 16 ---
$ws.Range("L16").Value = 0

# --- Row 18 ---
$ws.Range("H18").Value = 290.56
$ws.Range("L18").Value = 0

# --- Row 21 ---
$ws.Range("D21").Value = 1628.16
$ws.Range("G21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 13936.67

# --- Row 22 (summary "x de 20" counters) ---
$ws.Range("C22").Value = "1 de 20"
$ws.Range("D22").Value = "2 de 20"
$ws.Range("E22").Value = "2 de 20"
$ws.Range("F22").Value = "1 de 20"
$ws.Range("G22").Value = "0 de 20"
$ws.Range("H22").Value = "1 de 20"
$ws.Range("I22").Value = "2 de 20"
$ws.Range("J22").Value = "0 de 20"
$ws.Range("K22").Value = "2 de 20"
$ws.Range("L22").Value = "8 de 20"
$ws.Range("M22").Value = "1 de 20"
$ws.Range("N22").Value = "3 de 20"
